$wb = $excel.ActiveWorkbook

# New timestamp (Excel serial date) used for the appended row 52 on every sheet
$newTimestamp = [double]"45838.49256944445"

# Data for the new row 52 on each of the 4 worksheets.
# Column A (time) is the new timestamp above; columns B-I are copied
# unchanged from the previous last row (row 51) of each sheet.
$rowsData = @(
    @{
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x58"
        E = "0xf"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 344
        I = 15
    },
    @{
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x68"
        E = "0xe"
        F = 400
        G = [double]"5.68432987514711e+23"
        H = 360
        I = 14
    },
    @{
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x69"
        E = "0x3"
        F = 110
        G = [double]"5.68631262647114e+23"
        H = 105
        I = 3
    },
    @{
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x68"
        E = "0x3"
        F = 110
        G = [double]"9.85046333984776e+23"
        H = 104
        I = 3
    }
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $data = $rowsData[$i - 1]

    # Column A: new timestamp, matching the date style used by the rows above it
    $ws.Range("A52").NumberFormat = $ws.Range("A51").NumberFormat
    $ws.Range("A52").Value = $newTimestamp

    $ws.Range("B52").Value = $data.B
    $ws.Range("C52").Value = $data.C
    $ws.Range("D52").Value = $data.D
    $ws.Range("E52").Value = $data.E
    $ws.Range("F52").Value = $data.F
    $ws.Range("G52").Value = $data.G
    $ws.Range("H52").Value = $data.H
    $ws.Range("I52").Value = $data.I
}
